$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 176, shifting existing rows 176-233 down to 177-234.
$ws.Rows.Item(176).Insert()

# Populate the newly inserted row 176 with a new price record
# (same market/product template as the row that used to be here, with
# updated date, quality/volume/price and origin).
$ws.Cells.Item(176, 1).Value = 4
$ws.Cells.Item(176, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(176, 3).Value = "Los Lagos"
$ws.Cells.Item(176, 4).Value = 44795
$ws.Cells.Item(176, 5).Value = 10
$ws.Cells.Item(176, 6).Value = "Fruta"
$ws.Cells.Item(176, 7).Value = 100108
$ws.Cells.Item(176, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(176, 9).Value = 100108002
$ws.Cells.Item(176, 10).Value = "Mango"
$ws.Cells.Item(176, 11).Value = "Sin especificar"
$ws.Cells.Item(176, 12).Value = "Primera"
$ws.Cells.Item(176, 13).Value = 200
$ws.Cells.Item(176, 14).Value = 13000
$ws.Cells.Item(176, 15).Value = 14000
$ws.Cells.Item(176, 16).Value = 13500
$ws.Cells.Item(176, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(176, 18).Value = "Brasil"
$ws.Cells.Item(176, 19).Value = 3375
$ws.Cells.Item(176, 20).Value = 4
